# "Get presentation out of master slide view"
#
# The author had been poking around in Slide Master / Notes Master view and,
# on returning to Normal view and saving, PowerPoint re-typed/re-committed
# the "Acknowledgement" slide's text. The net effect on the two paragraphs
# below is that two adjacent runs that already shared identical run
# properties collapsed back into a single run (no wording changed).
#
# Slide 2, shape "Content Placeholder 2":
#   - "This material is based on work supported by the National Science
#      Foundation under Grants " + "DUE-1225708, ... DUE-1524877"
#        -> merges into one run
#   - "Copyright " + "and Licensing"
#        -> merges into one run

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(2)
$sh = $s.Shapes.Item(2)
Write-Host "Shape:" $sh.Name

$tf = $sh.TextFrame
$tr = $tf.TextRange

# --- Merge the "This material is based on ... Grants" + "DUE-...877" runs ---
$fullText = $tr.Text
$start1 = $fullText.IndexOf("This material is based on work supported") + 1
$len1 = "This material is based on work supported by the National Science Foundation under Grants DUE-1225708, DUE-1225738, DUE-1225688, DUE-1525039 DUE-1524898, and DUE-1524877".Length
$run1 = $tr.Characters($start1, $len1)
$run1.Text = "This material is based on work supported by the National Science Foundation under Grants DUE-1225708, DUE-1225738, DUE-1225688, DUE-1525039 DUE-1524898, and DUE-1524877"

# --- Merge the "Copyright " + "and Licensing" runs ---
$fullText2 = $tr.Text
$start2 = $fullText2.IndexOf("Copyright ") + 1
$len2 = "Copyright and Licensing".Length
$run2 = $tr.Characters($start2, $len2)
$run2.Text = "Copyright and Licensing"

Write-Host "Done."
